$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.765.23"
$ws.Range("E2").Value = "  +1.60%  "

$ws.Range("D3").Value = "2.265.63"
$ws.Range("E3").Value = "  +0.94%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "

$ws.Range("E7").Value = "  +2.16%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.31%  "

$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("E14").Value = "  +1.35%  "

$ws.Range("D15").Value = "2.616.21"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.33%  "

$ws.Range("D17").Value = "2.288.69"
$ws.Range("E17").Value = "  +6.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.770"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.85%  "

$ws.Range("D19").Value = "41.653.35"
$ws.Range("E19").Value = "  +1.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.68%  "

$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("E27").Value = "  +4.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.08%  "

$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("E30").Value = "  -0.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.13%  "

$ws.Range("E33").Value = "  +1.93%  "

$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("E35").Value = "  +1.77%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.01%  "

$ws.Range("E38").Value = "  +0.44%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.116"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.30%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.02%  "

$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("D43").Value = "2.013.32"
$ws.Range("E43").Value = "  -2.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.34%  "

$ws.Range("E45").Value = "  +1.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.82%  "

$ws.Range("E47").Value = "  +4.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.79%  "

$ws.Range("E50").Value = "  +1.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.21%  "
